# Apply the "Add data for 2022-10-01" update:
# - Rename sheet "Through 2022-09-22" -> "Through 2022-09-23"
# - Update header label "2022 (through 09-22)" -> "2022 (through 09-23)" (cell I1)
# - Update I10: 107 -> 111
# - Update I14 (Total): 1242 -> 1246

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-09-23"

$ws.Range("I1").Value = "2022 (through 09-23)"

$ws.Range("I10").Value = 111

$ws.Range("I14").Value = 1246
